$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Ensure LAST UPDATE column keeps storing plain text (not auto-converted to a date serial)
$ws.Range("I3:I5").NumberFormat = "@"

# Row 3
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = "04-Nov-2025"

# Row 4
$ws.Range("H4").Value = -99
$ws.Range("I4").Value = "04-Nov-2025"

# Row 5
$ws.Range("H5").Value = 286
$ws.Range("I5").Value = "04-Nov-2025"
